# Append two new registration rows to the bottom of the sheet (rows 12-13),
# extending the used range from A1:B11 to A1:B13.
#
# Column A holds account codes that look numeric ("08000", "12345") but must
# stay as literal text (leading zeros matter), so we prefix with a leading
# apostrophe to force text entry, exactly as a user typing into Excel would
# need to for a "General" formatted cell. Column B's timestamp strings are
# not valid numbers, so they are stored as text automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "'08000"
$ws.Range("B12").Value = "22/08/2025 18:50:41"

$ws.Range("A13").Value = "'12345"
$ws.Range("B13").Value = "22/08/2025 19:09:06"
